$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, New nombre_aides (C), New montant_total (E)
$updates = @(
    @{ Row = 8;   C = 1052;   E = 91427330 },
    @{ Row = 14;  C = 110819; E = 253257597 },
    @{ Row = 25;  C = 85736;  E = 374597636 },
    @{ Row = 91;  C = 151207; E = 482961698 },
    @{ Row = 92;  C = 409326; E = 1597623538 },
    @{ Row = 93;  C = 209674; E = 1310292697 },
    @{ Row = 94;  C = 94244;  E = 919337557 },
    @{ Row = 95;  C = 50811;  E = 934799343 },
    @{ Row = 96;  C = 17328;  E = 797814283 },
    @{ Row = 142; C = 168978; E = 681802746 },
    @{ Row = 148; C = 94;     E = 10649586 },
    @{ Row = 175; C = 80787;  E = 486200750 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
